$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the duplicate/erroneous row 8 (PERMNS_007 row with bad timestamp/env),
#    which shifts all subsequent rows up by one.
$ws.Rows.Item(8).Delete()

# 2. Insert 14 new blank rows starting at row 36 (pushes the former row 36
#    "PERMNS_094..." block down to rows 50-51).
$ws.Range("A36:A49").EntireRow.Insert()

# 3. Populate the newly inserted rows with the new test-case data.
$ws.Range("A36").Value = "PERMNS_035_Verify_Pricing_Permission_As_None"
$ws.Range("B36").Value = "Sorry, you do not have permissions to access this page."
$ws.Range("C36").Value = "Top displayed text is "
$ws.Range("D36").Value = "Permissions"
$ws.Range("E36").Value = "Passed"
$ws.Range("F36").Value = "2023-10-11 16:43:40"
$ws.Range("G36").Value = "QA Instance"

$ws.Range("A37").Value = "PERMNS_036_Verify_Discount Codes_Permission_As_None"
$ws.Range("B37").Value = "Sorry, you do not have permissions to access this page."
$ws.Range("C37").Value = "Top displayed text is "
$ws.Range("D37").Value = "Permissions"
$ws.Range("E37").Value = "Passed"
$ws.Range("F37").Value = "2023-10-11 16:44:13"
$ws.Range("G37").Value = "QA Instance"

$ws.Range("A38").Value = "PERMNS_037_Verify_Non Standard Pricing_Permission_As_None"
$ws.Range("B38").Value = "Sorry, you do not have permissions to access this page."
$ws.Range("C38").Value = "Top displayed text is "
$ws.Range("D38").Value = "Permissions"
$ws.Range("E38").Value = "Passed"
$ws.Range("F38").Value = "2023-10-11 16:44:46"
$ws.Range("G38").Value = "QA Instance"

$ws.Range("A39").Value = "PERMNS_038_Verify_Export_Permission_As_Yes_In_Pricing"
$ws.Range("B39").Value = "Centrifuge`nFilters`nAdd"
$ws.Range("C39").Value = "Top displayed text is Centrifuge`nFilters`nAdd"
$ws.Range("D39").Value = "Permissions"
$ws.Range("E39").Value = "Failed"
$ws.Range("F39").Value = "2023-10-11 16:45:22"
$ws.Range("G39").Value = "QA Instance"

$ws.Range("A40").Value = "PERMNS_039_Verify_Export_Permission_As_NoIn_Pricing"
$ws.Range("B40").Value = "Centrifuge`nExport`nFilters`nAdd"
$ws.Range("C40").Value = "Top displayed text is Centrifuge`nExport`nFilters`nAdd"
$ws.Range("D40").Value = "Permissions"
$ws.Range("E40").Value = "Passed"
$ws.Range("F40").Value = "2023-10-11 16:45:58"
$ws.Range("G40").Value = "QA Instance"

$ws.Range("A41").Value = "PERMNS_040_Verify_Import_Permission_As_Yes_In_Pricing"
$ws.Range("B41").Value = "Centrifuge`nExport`nFilters`nAdd"
$ws.Range("C41").Value = "Top displayed text is Centrifuge`nExport`nFilters`nAdd"
$ws.Range("D41").Value = "Permissions"
$ws.Range("E41").Value = "Failed"
$ws.Range("F41").Value = "2023-10-11 16:46:35"
$ws.Range("G41").Value = "QA Instance"

$ws.Range("A42").Value = "PERMNS_041_Verify_Import_Permission_As_No_In_Pricing"
$ws.Range("B42").Value = "Centrifuge`nImport`nExport`nFilters`nAdd"
$ws.Range("C42").Value = "Top displayed text is Centrifuge`nImport`nExport`nFilters`nAdd"
$ws.Range("D42").Value = "Permissions"
$ws.Range("E42").Value = "Passed"
$ws.Range("F42").Value = "2023-10-11 16:47:11"
$ws.Range("G42").Value = "QA Instance"

$ws.Range("A43").Value = "PERMNS_042_Verify_Export_Permission_As_Yes_In_Discount Codes"
$ws.Range("B43").Value = "Centrifuge`nMulti Edit`nFilters`nAdd"
$ws.Range("C43").Value = "Top displayed text is Centrifuge`nMulti Edit`nFilters`nAdd"
$ws.Range("D43").Value = "Permissions"
$ws.Range("E43").Value = "Failed"
$ws.Range("F43").Value = "2023-10-11 16:47:48"
$ws.Range("G43").Value = "QA Instance"

$ws.Range("A44").Value = "PERMNS_043_Verify_Export_Permission_As_No_In_Discount Codes"
$ws.Range("B44").Value = "Centrifuge`nMulti Edit`nExport`nFilters`nAdd"
$ws.Range("C44").Value = "Top displayed text is Centrifuge`nMulti Edit`nExport`nFilters`nAdd"
$ws.Range("D44").Value = "Permissions"
$ws.Range("E44").Value = "Passed"
$ws.Range("F44").Value = "2023-10-11 16:48:24"
$ws.Range("G44").Value = "QA Instance"

$ws.Range("A45").Value = "PERMNS_044_Verify_Export_Permission_As_Yes_In_Non Standard Pricing"
$ws.Range("B45").Value = "Pricing Rule Configurator`nExport"
$ws.Range("C45").Value = "Top displayed text is Pricing Rule Configurator`nExport"
$ws.Range("D45").Value = "Permissions"
$ws.Range("E45").Value = "Passed"
$ws.Range("F45").Value = "2023-10-11 16:49:08"
$ws.Range("G45").Value = "QA Instance"

$ws.Range("A46").Value = "PERMNS_045_Verify_Export_Permission_As_No_In_Non Standard Pricing"
$ws.Range("B46").Value = "Pricing Rule Configurator"
$ws.Range("C46").Value = "Top displayed text is Pricing Rule Configurator"
$ws.Range("D46").Value = "Permissions"
$ws.Range("E46").Value = "Passed"
$ws.Range("F46").Value = "2023-10-11 16:49:51"
$ws.Range("G46").Value = "QA Instance"

$ws.Range("A47").Value = "PERMNS_046_Verify_Pricing_Permission_As_View"
$ws.Range("B47").Value = "Centrifuge`nImport`nExport`nFilters`nAdd"
$ws.Range("C47").Value = "Top displayed text is Centrifuge`nImport`nExport`nFilters`nAdd"
$ws.Range("D47").Value = "Permissions"
$ws.Range("E47").Value = "Failed"
$ws.Range("F47").Value = "2023-10-11 16:50:25"
$ws.Range("G47").Value = "QA Instance"

$ws.Range("A48").Value = "PERMNS_047_Verify_Discount Codes_Permission_As_View"
$ws.Range("B48").Value = "Centrifuge`nMulti Edit`nExport`nFilters`nAdd"
$ws.Range("C48").Value = "Top displayed text is Centrifuge`nMulti Edit`nExport`nFilters`nAdd"
$ws.Range("D48").Value = "Permissions"
$ws.Range("E48").Value = "Failed"
$ws.Range("F48").Value = "2023-10-11 16:50:58"
$ws.Range("G48").Value = "QA Instance"

$ws.Range("A49").Value = "PERMNS_048_Verify_Non Standard Pricing_Permission_As_View"
$ws.Range("B49").Value = ""
$ws.Range("C49").Value = "Top displayed text is "
$ws.Range("D49").Value = "Permissions"
$ws.Range("E49").Value = "Passed"
$ws.Range("F49").Value = "2023-10-11 16:51:30"
$ws.Range("G49").Value = "QA Instance"

